$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 885094718.96
$ws.Range("P2").Value = 39603061814.54
$ws.Range("Q2").Value = 38381291885.18
$ws.Range("R2").Value = 15.6298165703
$ws.Range("S2").Value = 36093686798.35
$ws.Range("T2").Value = 36093686798.35
$ws.Range("U2").Value = 16.6783348605
$ws.Range("V2").Value = 177469349.34
$ws.Range("W2").Value = 975283703.36
$ws.Range("X2").Value = 801393169.6900001
$ws.Range("Y2").Value = 1345150216.77
$ws.Range("Z2").Value = 1306825118.95
$ws.Range("AA2").Value = 179903285.19
$ws.Range("AG2").Value = 198541219.42
$ws.Range("AP2").Value = 15.0857641909
$ws.Range("AQ2").Value = 3.146449377173
$ws.Range("AR2").Value = 24.812880062357
$ws.Range("AS2").Value = 738454661.61
$ws.Range("AT2").Value = -19.618511269903
